# "Generate Report for Handback"
# f165b98e-... finished its handback cycle (status flips from "Ready for
# handoff" to "Handed back: in sync with en-US", and gains Latest Target
# File / Latest Handback File / Latest Handback DateTime). Because the
# report sorts "done" rows above still-pending ones, f165b98e's row now
# sits above ce0e67cd's (still "Handback transform failed") on every
# sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview": rows 5/6 swap (f165b98e now above ce0e67cd); row 7
# (.localization-config) is unaffected.
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Cells.Item(5,1).Value = "f165b98e-5663-447e-8f1b-b1fc55318954.md"
$ov.Cells.Item(5,2).Value = "Handed back: in sync with en-US"
$ov.Cells.Item(5,3).Value = "Handed back: in sync with en-US"

$ov.Cells.Item(6,1).Value = "ce0e67cd-c4d9-428d-852e-f014ae0eae58.md"
$ov.Cells.Item(6,2).Value = "Handback transform failed"
$ov.Cells.Item(6,3).Value = "Handback transform failed"

$ovLinks = @($ov.Hyperlinks)
foreach ($hl in $ovLinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$5') {
        $hl.Address = "https://github.com/OpenLocalizationTest/oltest/blob/85befd01c027cb88b0d40e6a133f21f3af37d061/e2e/f165b98e-5663-447e-8f1b-b1fc55318954.md"
        $hl.TextToDisplay = "f165b98e-5663-447e-8f1b-b1fc55318954.md"
    } elseif ($addr -eq '$A$6') {
        $hl.Address = "https://github.com/OpenLocalizationTest/oltest/blob/c59ce95be3482f6698061bc08c60a6cb4c58ffe4/e2e/ce0e67cd-c4d9-428d-852e-f014ae0eae58.md"
        $hl.TextToDisplay = "ce0e67cd-c4d9-428d-852e-f014ae0eae58.md"
    }
}

# ---------------------------------------------------------------------
# Helper data + body shared by the "zh-cn" and "de-de" detail sheets.
# ---------------------------------------------------------------------
$langs = @("zh-cn", "de-de")

$handoffFile = @{
    "zh-cn" = "f165b98e-5663-447e-8f1b-b1fc55318954.3a5e58788516a282bb072181e7823ad9191fec77.zh-cn.xlf"
    "de-de" = "f165b98e-5663-447e-8f1b-b1fc55318954.3a5e58788516a282bb072181e7823ad9191fec77.de-de.xlf"
}
$handoffDate = @{
    "zh-cn" = "2016-03-09 14:20:18"
    "de-de" = "2016-03-09 14:20:20"
}
$handbackDate = @{
    "zh-cn" = "2016-03-09 14:20:56"
    "de-de" = "2016-03-09 14:21:01"
}
$handoffMdAddr = @{
    "zh-cn" = "https://github.com/OpenLocalizationTest/oltest/blob/85befd01c027cb88b0d40e6a133f21f3af37d061/e2e/f165b98e-5663-447e-8f1b-b1fc55318954.md"
    "de-de" = "https://github.com/OpenLocalizationTest/oltest/blob/85befd01c027cb88b0d40e6a133f21f3af37d061/e2e/f165b98e-5663-447e-8f1b-b1fc55318954.md"
}
$handoffXlfAddr = @{
    "zh-cn" = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/970bf7c6a49f09c54d11624786483d1f543453f1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/f165b98e-5663-447e-8f1b-b1fc55318954.3a5e58788516a282bb072181e7823ad9191fec77.zh-cn.xlf"
    "de-de" = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/197452497f1d42b0babc165f9fe201be1cf9d994/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/f165b98e-5663-447e-8f1b-b1fc55318954.3a5e58788516a282bb072181e7823ad9191fec77.de-de.xlf"
}
$targetMdAddr = @{
    "zh-cn" = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/e6d6a5d6be93dbbef64db8db2f1a22404ab8f0f9/e2e/f165b98e-5663-447e-8f1b-b1fc55318954.md"
    "de-de" = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/0b1a2e7e2f19cb2c56ea5f3b1f4a55df4f2f6a3e/e2e/f165b98e-5663-447e-8f1b-b1fc55318954.md"
}
$handbackXlfAddr = @{
    "zh-cn" = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/3e7c9e0a6a4b2f8d5c2e7f1b9a6d4c8e2f0b5a7d/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/f165b98e-5663-447e-8f1b-b1fc55318954.3a5e58788516a282bb072181e7823ad9191fec77.zh-cn.xlf"
    "de-de" = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/9a4d2b6e8f1c3a5d7b0e9f2c4a6d8b1e3f5a7c9d/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/f165b98e-5663-447e-8f1b-b1fc55318954.3a5e58788516a282bb072181e7823ad9191fec77.de-de.xlf"
}

$ceFile = @{
    "zh-cn" = "ce0e67cd-c4d9-428d-852e-f014ae0eae58.e8986e400d09cce0c6c2d335f4cdabfa2645f548.zh-cn.xlf"
    "de-de" = "ce0e67cd-c4d9-428d-852e-f014ae0eae58.e8986e400d09cce0c6c2d335f4cdabfa2645f548.de-de.xlf"
}
$ceDate = @{
    "zh-cn" = "2016-03-09 14:17:31"
    "de-de" = "2016-03-09 14:17:35"
}
$ceMdAddr = @{
    "zh-cn" = "https://github.com/OpenLocalizationTest/oltest/blob/c59ce95be3482f6698061bc08c60a6cb4c58ffe4/e2e/ce0e67cd-c4d9-428d-852e-f014ae0eae58.md"
    "de-de" = "https://github.com/OpenLocalizationTest/oltest/blob/c59ce95be3482f6698061bc08c60a6cb4c58ffe4/e2e/ce0e67cd-c4d9-428d-852e-f014ae0eae58.md"
}
$ceXlfAddr = @{
    "zh-cn" = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6bbee925604a3a023356de7a2111b12db76dc5bc/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/ce0e67cd-c4d9-428d-852e-f014ae0eae58.e8986e400d09cce0c6c2d335f4cdabfa2645f548.zh-cn.xlf"
    "de-de" = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c0115fc545beaef72ac83cb879c2b1f87f3da87f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/ce0e67cd-c4d9-428d-852e-f014ae0eae58.e8986e400d09cce0c6c2d335f4cdabfa2645f548.de-de.xlf"
}

foreach ($lang in $langs) {
    $ws = $wb.Worksheets.Item($lang)

    # --- row 5 becomes f165b98e's completed handback entry ---
    $ws.Cells.Item(5,1).Value = "f165b98e-5663-447e-8f1b-b1fc55318954.md"
    $ws.Cells.Item(5,2).Value = "Handed back: in sync with en-US"
    $ws.Cells.Item(5,3).Value = $handoffFile[$lang]
    $ws.Cells.Item(5,4).Value = $handoffDate[$lang]
    $ws.Cells.Item(5,5).Value = "f165b98e-5663-447e-8f1b-b1fc55318954.md"
    $ws.Cells.Item(5,6).Value = $handoffFile[$lang]
    $ws.Cells.Item(5,7).Value = $handbackDate[$lang]
    $ws.Cells.Item(5,8).Value = "Include"

    # --- row 6 becomes ce0e67cd's (still-failed) entry ---
    $ws.Cells.Item(6,1).Value = "ce0e67cd-c4d9-428d-852e-f014ae0eae58.md"
    $ws.Cells.Item(6,2).Value = "Handback transform failed"
    $ws.Cells.Item(6,3).Value = $ceFile[$lang]
    $ws.Cells.Item(6,4).Value = $ceDate[$lang]
    $ws.Cells.Item(6,7).Value = "0001-01-01 00:00:00"
    $ws.Cells.Item(6,8).Value = "Include"

    # --- existing hyperlinks on A5/C5/A6/C6 get repointed to the new
    #     row owners, and two new ones appear on E5/F5 ---
    $links = @($ws.Hyperlinks)
    foreach ($hl in $links) {
        $addr = $hl.Range.Address()
        if ($addr -eq '$A$5') {
            $hl.Address = $handoffMdAddr[$lang]
            $hl.TextToDisplay = "f165b98e-5663-447e-8f1b-b1fc55318954.md"
        } elseif ($addr -eq '$C$5') {
            $hl.Address = $handoffXlfAddr[$lang]
            $hl.TextToDisplay = $handoffFile[$lang]
        } elseif ($addr -eq '$A$6') {
            $hl.Address = $ceMdAddr[$lang]
            $hl.TextToDisplay = "ce0e67cd-c4d9-428d-852e-f014ae0eae58.md"
        } elseif ($addr -eq '$C$6') {
            $hl.Address = $ceXlfAddr[$lang]
            $hl.TextToDisplay = $ceFile[$lang]
        }
    }

    $ws.Hyperlinks.Add($ws.Range("E5"), $targetMdAddr[$lang], "", "", "f165b98e-5663-447e-8f1b-b1fc55318954.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("F5"), $handbackXlfAddr[$lang], "", "", $handoffFile[$lang]) | Out-Null
}
